$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1299
$ws.Range("I15").Value = 1299
$ws.Range("K15").Value = 3897
$ws.Range("M15").Value = -3728
$ws.Range("H28").Value = 634.8387
$ws.Range("I28").Value = 244
$ws.Range("J28").Value = 1176
$ws.Range("K28").Value = 244
$ws.Range("L28").Value = 1176
$ws.Range("M28").Value = 241
$ws.Range("N28").Value = -2146
$ws.Range("H82").Value = 2468.6667
$ws.Range("I82").Value = 784.2
$ws.Range("K82").Value = 2352.6
$ws.Range("M82").Value = -1946.6
$ws.Range("H85").Value = 2468.6667
$ws.Range("I85").Value = 784.2
$ws.Range("K85").Value = 2352.6
$ws.Range("M85").Value = -948.6000000000004
$ws.Range("H98").Value = 2212
$ws.Range("I98").Value = 2098.5715
$ws.Range("J98").Value = 3006
$ws.Range("K98").Value = 2098.5715
$ws.Range("L98").Value = 3006
$ws.Range("M98").Value = -600.5715
$ws.Range("N98").Value = -6002
$ws.Range("H122").Value = 2212
$ws.Range("I122").Value = 2098.5715
$ws.Range("J122").Value = 3006
$ws.Range("K122").Value = 6295.7145
$ws.Range("L122").Value = 9018
$ws.Range("M122").Value = -3845.7145
$ws.Range("N122").Value = -13918

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22549.568
$ws.Range("I32").Value = 3720.6428
$ws.Range("J32").Value = 142370
$ws.Range("K32").Value = 3720.6428
$ws.Range("L32").Value = 142370
$ws.Range("M32").Value = -3433.6428
$ws.Range("N32").Value = -142944
$ws.Range("H74").Value = 4547564.5
$ws.Range("I74").Value = 1461.6
$ws.Range("J74").Value = 14289214
$ws.Range("K74").Value = 1461.6
$ws.Range("L74").Value = 14289214
$ws.Range("M74").Value = -587.5999999999999
$ws.Range("N74").Value = -14290962
$ws.Range("H77").Value = 4547564.5
$ws.Range("I77").Value = 1461.6
$ws.Range("J77").Value = 14289214
$ws.Range("K77").Value = 7308
$ws.Range("L77").Value = 71446070
$ws.Range("M77").Value = -2940
$ws.Range("N77").Value = -71454806
$ws.Range("H133").Value = 45000
$ws.Range("J133").Value = 45000
$ws.Range("L133").Value = 45000
$ws.Range("N133").Value = -50060

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 27416.666
$ws.Range("I8").Value = 950
$ws.Range("J8").Value = 40650
$ws.Range("K8").Value = 950
$ws.Range("L8").Value = 40650
$ws.Range("M8").Value = -810
$ws.Range("N8").Value = -40930

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H88").Value = 47975
$ws.Range("J88").Value = 47975
$ws.Range("L88").Value = 47975
$ws.Range("N88").Value = -48787
$ws.Range("H91").Value = 47975
$ws.Range("J91").Value = 47975
$ws.Range("L91").Value = 47975
$ws.Range("N91").Value = -50783
$ws.Range("H99").Value = 10600.733
$ws.Range("I99").Value = 3343.2856
$ws.Range("J99").Value = 16951
$ws.Range("K99").Value = 3343.2856
$ws.Range("L99").Value = 16951
$ws.Range("M99").Value = -1845.2856
$ws.Range("N99").Value = -19947
$ws.Range("H122").Value = 448.4
$ws.Range("I122").Value = 423.25
$ws.Range("J122").Value = 549
$ws.Range("K122").Value = 1269.75
$ws.Range("L122").Value = 1647
$ws.Range("M122").Value = 1180.25
$ws.Range("N122").Value = -6547
$ws.Range("H126").Value = 10600.733
$ws.Range("I126").Value = 3343.2856
$ws.Range("J126").Value = 16951
$ws.Range("K126").Value = 10029.8568
$ws.Range("L126").Value = 50853
$ws.Range("M126").Value = -7559.856800000001
$ws.Range("N126").Value = -55793
$ws.Range("H133").Value = 50995.6
$ws.Range("J133").Value = 50995.6
$ws.Range("L133").Value = 50995.6
$ws.Range("N133").Value = -56055.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 622.75
$ws.Range("I63").Value = 459
$ws.Range("J63").Value = 1114
$ws.Range("K63").Value = 1377
$ws.Range("L63").Value = 3342
$ws.Range("M63").Value = -628
$ws.Range("N63").Value = -4840
$ws.Range("H66").Value = 622.75
$ws.Range("I66").Value = 459
$ws.Range("J66").Value = 1114
$ws.Range("K66").Value = 4131
$ws.Range("L66").Value = 10026
$ws.Range("M66").Value = -387
$ws.Range("N66").Value = -17514
$ws.Range("H131").Value = 785.5700000000001
$ws.Range("I131").Value = 404.69232
$ws.Range("J131").Value = 842.4828
$ws.Range("K131").Value = 1214.07696
$ws.Range("L131").Value = 2527.4484
$ws.Range("M131").Value = 3825.92304
$ws.Range("N131").Value = -12607.4484

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H102").Value = 2367.5
$ws.Range("J102").Value = 2813.5
$ws.Range("L102").Value = 2813.5
$ws.Range("N102").Value = -6057.5
$ws.Range("H117").Value = 19999
$ws.Range("J117").Value = 19999
$ws.Range("L117").Value = 19999
$ws.Range("N117").Value = -26883
$ws.Range("H122").Value = 1883
$ws.Range("I122").Value = 1949.5714
$ws.Range("J122").Value = 1650
$ws.Range("K122").Value = 5848.7142
$ws.Range("L122").Value = 4950
$ws.Range("M122").Value = -3398.7142
$ws.Range("N122").Value = -9850
$ws.Range("H135").Value = 28725.578
$ws.Range("J135").Value = 28725.578
$ws.Range("L135").Value = 28725.578
$ws.Range("N135").Value = -38865.578

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4150
$ws.Range("I7").Value = 2000
$ws.Range("J7").Value = 5225
$ws.Range("K7").Value = 2000
$ws.Range("L7").Value = 5225
$ws.Range("M7").Value = -1888
$ws.Range("N7").Value = -5449
$ws.Range("H81").Value = 44181
$ws.Range("J81").Value = 44181
$ws.Range("L81").Value = 44181
$ws.Range("N81").Value = -46177
$ws.Range("H84").Value = 44181
$ws.Range("J84").Value = 44181
$ws.Range("L84").Value = 132543
$ws.Range("N84").Value = -142527
$ws.Range("H126").Value = 4150
$ws.Range("I126").Value = 2000
$ws.Range("J126").Value = 5225
$ws.Range("K126").Value = 6000
$ws.Range("L126").Value = 15675
$ws.Range("M126").Value = -3530
$ws.Range("N126").Value = -20615

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1878.25
$ws.Range("I122").Value = 1837.6666
$ws.Range("K122").Value = 5512.9998
$ws.Range("M122").Value = -3062.9998
$ws.Range("H126").Value = 1684.875
$ws.Range("I126").Value = 1496.7693
$ws.Range("J126").Value = 2500
$ws.Range("K126").Value = 4490.3079
$ws.Range("L126").Value = 7500
$ws.Range("M126").Value = -2020.3079
$ws.Range("N126").Value = -12440
